$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2221.1667
$ws.Range("I28").Value = 350.53333
$ws.Range("J28").Value = 5338.8887
$ws.Range("K28").Value = 350.53333
$ws.Range("L28").Value = 5338.8887
$ws.Range("M28").Value = 134.46667
$ws.Range("N28").Value = -6308.8887
$ws.Range("H47").Value = 15213.375
$ws.Range("I47").Value = 12772.167
$ws.Range("J47").Value = 22537
$ws.Range("K47").Value = 12772.167
$ws.Range("L47").Value = 22537
$ws.Range("M47").Value = -11800.167
$ws.Range("N47").Value = -24481
$ws.Range("H55").Value = 208.6875
$ws.Range("I55").Value = 133.16667
$ws.Range("J55").Value = 254
$ws.Range("K55").Value = 133.16667
$ws.Range("L55").Value = 254
$ws.Range("M55").Value = 80.83332999999999
$ws.Range("N55").Value = -682
$ws.Range("H69").Value = 4620
$ws.Range("I69").Value = 4633.3335
$ws.Range("K69").Value = 13900.0005
$ws.Range("M69").Value = -13026.0005
$ws.Range("H72").Value = 4620
$ws.Range("I72").Value = 4633.3335
$ws.Range("K72").Value = 41700.0015
$ws.Range("M72").Value = -37332.0015
$ws.Range("H80").Value = 2399.0688
$ws.Range("I80").Value = 2960.2942
$ws.Range("K80").Value = 8880.882599999999
$ws.Range("M80").Value = -7882.882599999999
$ws.Range("H83").Value = 2399.0688
$ws.Range("I83").Value = 2960.2942
$ws.Range("K83").Value = 26642.6478
$ws.Range("M83").Value = -21650.6478
$ws.Range("H86").Value = 3130.1304
$ws.Range("I86").Value = 2199.4119
$ws.Range("J86").Value = 5767.1665
$ws.Range("K86").Value = 2199.4119
$ws.Range("L86").Value = 5767.1665
$ws.Range("M86").Value = -1076.4119
$ws.Range("N86").Value = -8013.1665
$ws.Range("H88").Value = 1535
$ws.Range("I88").Value = 3024.25
$ws.Range("J88").Value = 542.1667
$ws.Range("K88").Value = 3024.25
$ws.Range("L88").Value = 542.1667
$ws.Range("M88").Value = -2618.25
$ws.Range("N88").Value = -1354.1667
$ws.Range("H89").Value = 3130.1304
$ws.Range("I89").Value = 2199.4119
$ws.Range("J89").Value = 5767.1665
$ws.Range("K89").Value = 10997.0595
$ws.Range("L89").Value = 28835.8325
$ws.Range("M89").Value = -5381.059499999999
$ws.Range("N89").Value = -40067.8325
$ws.Range("H91").Value = 1535
$ws.Range("I91").Value = 3024.25
$ws.Range("J91").Value = 542.1667
$ws.Range("K91").Value = 3024.25
$ws.Range("L91").Value = 542.1667
$ws.Range("M91").Value = -1620.25
$ws.Range("N91").Value = -3350.1667
$ws.Range("H107").Value = 1554.3684
$ws.Range("I107").Value = 1277
$ws.Range("J107").Value = 2331
$ws.Range("K107").Value = 1277
$ws.Range("L107").Value = 2331
$ws.Range("M107").Value = 643
$ws.Range("N107").Value = -6171
$ws.Range("H111").Value = 5634.2
$ws.Range("I111").Value = 5817.875
$ws.Range("K111").Value = 17453.625
$ws.Range("M111").Value = -14386.625
$ws.Range("H112").Value = 6099737
$ws.Range("J112").Value = 6412401.5
$ws.Range("L112").Value = 19237204.5
$ws.Range("N112").Value = -19239420.5
$ws.Range("H113").Value = 6100.385
$ws.Range("I113").Value = 4967.5
$ws.Range("J113").Value = 7071.4287
$ws.Range("K113").Value = 4967.5
$ws.Range("L113").Value = 7071.4287
$ws.Range("M113").Value = -1713.5
$ws.Range("N113").Value = -13579.4287
$ws.Range("H116").Value = 5167.222
$ws.Range("I116").Value = 5143.5713
$ws.Range("K116").Value = 5143.5713
$ws.Range("M116").Value = -1701.5713
$ws.Range("H118").Value = 3840.9167
$ws.Range("I118").Value = 1488.1111
$ws.Range("K118").Value = 4464.3333
$ws.Range("M118").Value = -2807.3333
$ws.Range("H137").Value = 39217540
$ws.Range("I137").Value = 22224096
$ws.Range("K137").Value = 66672288
$ws.Range("M137").Value = -66669738
$ws.Range("H141").Value = 2528.2
$ws.Range("I141").Value = 2651.4285
$ws.Range("J141").Value = 803
$ws.Range("K141").Value = 7954.2855
$ws.Range("L141").Value = 2409
$ws.Range("M141").Value = -2774.2855
$ws.Range("N141").Value = -12769

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2427.95
$ws.Range("I2").Value = 2658.7693
$ws.Range("J2").Value = 1999.2858
$ws.Range("K2").Value = 2658.7693
$ws.Range("L2").Value = 1999.2858
$ws.Range("M2").Value = -2545.7693
$ws.Range("N2").Value = -2225.2858
$ws.Range("H6").Value = 668500.7
$ws.Range("I6").Value = 2751
$ws.Range("K6").Value = 2751
$ws.Range("M6").Value = -2578
$ws.Range("H61").Value = 27032092
$ws.Range("I61").Value = 41670120
$ws.Range("J61").Value = 8044.615
$ws.Range("K61").Value = 41670120
$ws.Range("L61").Value = 8044.615
$ws.Range("M61").Value = -41669908
$ws.Range("N61").Value = -8468.615
$ws.Range("H63").Value = 2845
$ws.Range("I63").Value = 2845
$ws.Range("K63").Value = 2845
$ws.Range("M63").Value = -2159
$ws.Range("H66").Value = 2845
$ws.Range("I66").Value = 2845
$ws.Range("K66").Value = 14225
$ws.Range("M66").Value = -10793
$ws.Range("H74").Value = 29446158
$ws.Range("I74").Value = 29446158
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 29446158
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -29445284
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 29446158
$ws.Range("I77").Value = 29446158
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 147230790
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -147226422
$ws.Range("N77").ClearContents()
$ws.Range("H88").Value = 1913.8334
$ws.Range("J88").Value = 2661.3333
$ws.Range("L88").Value = 2661.3333
$ws.Range("N88").Value = -3473.3333
$ws.Range("H91").Value = 1913.8334
$ws.Range("J91").Value = 2661.3333
$ws.Range("L91").Value = 2661.3333
$ws.Range("N91").Value = -5469.3333
$ws.Range("H102").Value = 4204.25
$ws.Range("I102").Value = 3998.5
$ws.Range("K102").Value = 3998.5
$ws.Range("M102").Value = -2376.5
$ws.Range("H110").Value = 11566.723
$ws.Range("I110").Value = 12716.193
$ws.Range("K110").Value = 12716.193
$ws.Range("M110").Value = -10671.193
$ws.Range("H116").Value = 2427.95
$ws.Range("I116").Value = 2658.7693
$ws.Range("J116").Value = 1999.2858
$ws.Range("K116").Value = 2658.7693
$ws.Range("L116").Value = 1999.2858
$ws.Range("M116").Value = -364.7692999999999
$ws.Range("N116").Value = -6587.2858
$ws.Range("H122").Value = 2411.1177
$ws.Range("I122").Value = 1718.8948
$ws.Range("K122").Value = 5156.6844
$ws.Range("M122").Value = -2706.6844
$ws.Range("H136").Value = 27032092
$ws.Range("I136").Value = 41670120
$ws.Range("J136").Value = 8044.615
$ws.Range("K136").Value = 125010360
$ws.Range("L136").Value = 24133.845
$ws.Range("M136").Value = -125007810
$ws.Range("N136").Value = -29233.845

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2427.95
$ws.Range("I3").Value = 2658.7693
$ws.Range("J3").Value = 1999.2858
$ws.Range("K3").Value = 2658.7693
$ws.Range("L3").Value = 1999.2858
$ws.Range("M3").Value = -2544.7693
$ws.Range("N3").Value = -2227.2858
$ws.Range("H86").Value = 19435.77
$ws.Range("I86").Value = 3020.625
$ws.Range("K86").Value = 3020.625
$ws.Range("M86").Value = -1897.625
$ws.Range("H89").Value = 19435.77
$ws.Range("I89").Value = 3020.625
$ws.Range("K89").Value = 15103.125
$ws.Range("M89").Value = -9487.125
$ws.Range("H94").Value = 2182.9167
$ws.Range("I94").Value = 1774.5
$ws.Range("J94").Value = 2999.75
$ws.Range("K94").Value = 1774.5
$ws.Range("L94").Value = 2999.75
$ws.Range("M94").Value = -1323.5
$ws.Range("N94").Value = -3901.75
$ws.Range("H106").Value = 142725
$ws.Range("J106").Value = 142725
$ws.Range("L106").Value = 142725
$ws.Range("N106").Value = -145249

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26319334
$ws.Range("I31").Value = 3022.7666
$ws.Range("J31").Value = 125005496
$ws.Range("K31").Value = 3022.7666
$ws.Range("L31").Value = 125005496
$ws.Range("M31").Value = -2727.7666
$ws.Range("N31").Value = -125006086
$ws.Range("H34").Value = 26319334
$ws.Range("I34").Value = 3022.7666
$ws.Range("J34").Value = 125005496
$ws.Range("K34").Value = 3022.7666
$ws.Range("L34").Value = 125005496
$ws.Range("M34").Value = -2820.7666
$ws.Range("N34").Value = -125005900
$ws.Range("H62").Value = 4831.8237
$ws.Range("J62").Value = 6064.4287
$ws.Range("L62").Value = 6064.4287
$ws.Range("N62").Value = -7312.4287
$ws.Range("H65").Value = 4831.8237
$ws.Range("J65").Value = 6064.4287
$ws.Range("L65").Value = 30322.1435
$ws.Range("N65").Value = -36562.14350000001
$ws.Range("H95").Value = 390186.5
$ws.Range("J95").Value = 390186.5
$ws.Range("L95").Value = 390186.5
$ws.Range("N95").Value = -395678.5
$ws.Range("H99").Value = 6168.625
$ws.Range("I99").Value = 7450
$ws.Range("K99").Value = 7450
$ws.Range("M99").Value = -5952
$ws.Range("H107").Value = 1903.6111
$ws.Range("I107").Value = 1228.909
$ws.Range("K107").Value = 1228.909
$ws.Range("M107").Value = 691.0909999999999
$ws.Range("H126").Value = 6168.625
$ws.Range("I126").Value = 7450
$ws.Range("K126").Value = 22350
$ws.Range("M126").Value = -19880
$ws.Range("H132").Value = 1733.85
$ws.Range("I132").Value = 1357.9445
$ws.Range("K132").Value = 4073.8335
$ws.Range("M132").Value = -1543.8335

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1310.2222
$ws.Range("I5").Value = 536.35297
$ws.Range("J5").Value = 2625.8
$ws.Range("K5").Value = 1609.05891
$ws.Range("L5").Value = 7877.400000000001
$ws.Range("M5").Value = -1497.05891
$ws.Range("N5").Value = -8101.400000000001
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H55").Value = 10417233
$ws.Range("I55").Value = 16667023
$ws.Range("J55").Value = 6667359
$ws.Range("K55").Value = 50001069
$ws.Range("L55").Value = 20002077
$ws.Range("M55").Value = -50000892
$ws.Range("N55").Value = -20002431
$ws.Range("H86").Value = 804.9231
$ws.Range("I86").Value = 728.3333
$ws.Range("J86").Value = 977.25
$ws.Range("K86").Value = 2184.9999
$ws.Range("L86").Value = 2931.75
$ws.Range("M86").Value = -998.9998999999998
$ws.Range("N86").Value = -5303.75
$ws.Range("H89").Value = 804.9231
$ws.Range("I89").Value = 728.3333
$ws.Range("J89").Value = 977.25
$ws.Range("K89").Value = 6554.9997
$ws.Range("L89").Value = 8795.25
$ws.Range("M89").Value = -626.9997000000003
$ws.Range("N89").Value = -20651.25
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 400
$ws.Range("K92").Value = 1200
$ws.Range("M92").Value = 48
$ws.Range("H107").Value = 963
$ws.Range("I107").Value = 651.8
$ws.Range("J107").Value = 1352
$ws.Range("K107").Value = 1955.4
$ws.Range("L107").Value = 4056
$ws.Range("M107").Value = -35.39999999999986
$ws.Range("N107").Value = -7896
$ws.Range("H132").Value = 7413750.5
$ws.Range("J132").Value = 7413750.5
$ws.Range("L132").Value = 66723754.5
$ws.Range("N132").Value = -66728814.5
$ws.Range("H135").Value = 1310.2222
$ws.Range("I135").Value = 536.35297
$ws.Range("J135").Value = 2625.8
$ws.Range("K135").Value = 4827.17673
$ws.Range("L135").Value = 23632.2
$ws.Range("M135").Value = -2292.17673
$ws.Range("N135").Value = -28702.2
$ws.Range("H140").Value = 1608
$ws.Range("I140").Value = 1608
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4824
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 356
$ws.Range("N140").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H80").Value = 3742.8462
$ws.Range("I80").Value = 3289.8
$ws.Range("K80").Value = 3289.8
$ws.Range("M80").Value = -2291.8
$ws.Range("H83").Value = 3742.8462
$ws.Range("I83").Value = 3289.8
$ws.Range("K83").Value = 16449
$ws.Range("M83").Value = -11457
$ws.Range("H113").Value = 3227.762
$ws.Range("I113").Value = 2840.3333
$ws.Range("J113").Value = 3744.3333
$ws.Range("K113").Value = 2840.3333
$ws.Range("L113").Value = 3744.3333
$ws.Range("M113").Value = -670.3332999999998
$ws.Range("N113").Value = -8084.3333
$ws.Range("H126").Value = 12005292
$ws.Range("J126").Value = 20004344
$ws.Range("L126").Value = 60013032
$ws.Range("N126").Value = -60017972
$ws.Range("H132").Value = 1902.814
$ws.Range("I132").Value = 1639.9459
$ws.Range("K132").Value = 4919.8377
$ws.Range("M132").Value = -2389.8377

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3831.182
$ws.Range("I7").Value = 3508.5881
$ws.Range("K7").Value = 3508.5881
$ws.Range("M7").Value = -3396.5881
$ws.Range("H16").Value = 679.9048
$ws.Range("J16").Value = 1199
$ws.Range("L16").Value = 1199
$ws.Range("N16").Value = -1539
$ws.Range("H22").Value = 4088.6
$ws.Range("I22").Value = 1397.25
$ws.Range("J22").Value = 5882.8335
$ws.Range("K22").Value = 1397.25
$ws.Range("L22").Value = 5882.8335
$ws.Range("M22").Value = -1102.25
$ws.Range("N22").Value = -6472.8335
$ws.Range("H27").Value = 4088.6
$ws.Range("I27").Value = 1397.25
$ws.Range("J27").Value = 5882.8335
$ws.Range("K27").Value = 1397.25
$ws.Range("L27").Value = 5882.8335
$ws.Range("M27").Value = -1290.25
$ws.Range("N27").Value = -6096.8335
$ws.Range("H40").Value = 4217.7896
$ws.Range("I40").Value = 4008.2222
$ws.Range("K40").Value = 4008.2222
$ws.Range("M40").Value = -3872.2222
$ws.Range("H46").Value = 1933.8529
$ws.Range("I46").Value = 1066.619
$ws.Range("J46").Value = 3334.7693
$ws.Range("K46").Value = 1066.619
$ws.Range("L46").Value = 3334.7693
$ws.Range("M46").Value = -878.6189999999999
$ws.Range("N46").Value = -3710.7693
$ws.Range("H68").Value = 4070.4443
$ws.Range("I68").Value = 3655
$ws.Range("J68").Value = 5524.5
$ws.Range("K68").Value = 3655
$ws.Range("L68").Value = 5524.5
$ws.Range("M68").Value = -2906
$ws.Range("N68").Value = -7022.5
$ws.Range("H71").Value = 4070.4443
$ws.Range("I71").Value = 3655
$ws.Range("J71").Value = 5524.5
$ws.Range("K71").Value = 18275
$ws.Range("L71").Value = 27622.5
$ws.Range("M71").Value = -14531
$ws.Range("N71").Value = -35110.5
$ws.Range("H100").Value = 2825.5293
$ws.Range("I100").Value = 1998.7142
$ws.Range("J100").Value = 3404.3
$ws.Range("K100").Value = 1998.7142
$ws.Range("L100").Value = 3404.3
$ws.Range("M100").Value = -1457.7142
$ws.Range("N100").Value = -4486.3
$ws.Range("H102").Value = 59332.668
$ws.Range("I102").Value = 57499.5
$ws.Range("J102").Value = 62999
$ws.Range("K102").Value = 57499.5
$ws.Range("L102").Value = 62999
$ws.Range("M102").Value = -54254.5
$ws.Range("N102").Value = -69489
$ws.Range("H111").Value = 105193.5
$ws.Range("J111").Value = 105193.5
$ws.Range("L111").Value = 105193.5
$ws.Range("N111").Value = -113373.5
$ws.Range("H122").Value = 4050.2104
$ws.Range("I122").Value = 3585.5293
$ws.Range("K122").Value = 10756.5879
$ws.Range("M122").Value = -8306.5879
$ws.Range("H126").Value = 3831.182
$ws.Range("I126").Value = 3508.5881
$ws.Range("K126").Value = 10525.7643
$ws.Range("M126").Value = -8055.764299999999
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H136").Value = 3027.4473
$ws.Range("I136").Value = 3055.2163
$ws.Range("K136").Value = 9165.6489
$ws.Range("M136").Value = -6615.6489
$ws.Range("H138").Value = 89999
$ws.Range("J138").Value = 89999
$ws.Range("L138").Value = 89999
$ws.Range("N138").Value = -100279

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1072.2424
$ws.Range("I81").Value = 947.871
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 1895.742
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -834.742
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1072.2424
$ws.Range("I84").Value = 947.871
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 9478.709999999999
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -4174.709999999999
$ws.Range("N84").Value = -40608
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 686.5833
$ws.Range("I107").Value = 626.8889
$ws.Range("J107").Value = 865.6667
$ws.Range("K107").Value = 1880.6667
$ws.Range("L107").Value = 2597.0001
$ws.Range("M107").Value = 39.33329999999978
$ws.Range("N107").Value = -6437.0001
$ws.Range("H113").Value = 866.2
$ws.Range("I113").Value = 409.14285
$ws.Range("K113").Value = 1227.42855
$ws.Range("M113").Value = 942.5714499999999
$ws.Range("H114").Value = 79998
$ws.Range("J114").Value = 79998
$ws.Range("L114").Value = 79998
$ws.Range("N114").Value = -88676
$ws.Range("H126").Value = 4241.92
$ws.Range("I126").Value = 4348.0835
$ws.Range("J126").Value = 1694
$ws.Range("K126").Value = 13044.2505
$ws.Range("L126").Value = 5082
$ws.Range("M126").Value = -10574.2505
$ws.Range("N126").Value = -10022
$ws.Range("H128").Value = 35325
$ws.Range("I128").Value = 20650
$ws.Range("J128").Value = 50000
$ws.Range("K128").Value = 20650
$ws.Range("L128").Value = 50000
$ws.Range("M128").Value = -15670
$ws.Range("N128").Value = -59960
$ws.Range("H132").Value = 3355.4177
$ws.Range("I132").Value = 3574.9385
$ws.Range("J132").Value = 2336.2144
$ws.Range("K132").Value = 10724.8155
$ws.Range("L132").Value = 7008.6432
$ws.Range("M132").Value = -8194.815500000001
$ws.Range("N132").Value = -12068.6432
$ws.Range("H136").Value = 1925.4524
$ws.Range("I136").Value = 1862.6586
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 5587.9758
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -3037.9758
$ws.Range("N136").Value = -18600
